$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23 (shifts existing rows 23:56 down to 24:57)
$ws.Range("A23").EntireRow.Insert()

# Populate the newly inserted row 23 with the new weekly price observation
$ws.Range("A23").Value = 5
$ws.Range("B23").Value = "Macroferia Regional de Talca"
$ws.Range("C23").Value = "Maule"
$ws.Range("D23").Value = 44587
$ws.Range("E23").Value = 7
$ws.Range("F23").Value = "Fruta"
$ws.Range("G23").Value = 100101
$ws.Range("H23").Value = "Berries"
$ws.Range("I23").Value = 100101001
$ws.Range("J23").Value = "Arándano (blue)"
$ws.Range("K23").Value = "Sin especificar"
$ws.Range("L23").Value = "Primera"
$ws.Range("M23").Value = 150
$ws.Range("N23").Value = 3500
$ws.Range("O23").Value = 3500
$ws.Range("P23").Value = 3500
$ws.Range("Q23").Value = "$/bandeja 2 kilos"
$ws.Range("R23").Value = "Provincia de Linares"
$ws.Range("S23").Value = 1750
$ws.Range("T23").Value = 2
